$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.803.24"
$ws.Range("E2").Value = "'  +1.92%  "
$ws.Range("D3").Value = "'1.657.39"
$ws.Range("E3").Value = "'  +1.88%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "'  +0.10%  "
$ws.Range("D6").Value = "'304.17"
$ws.Range("E6").Value = "'  +0.57%  "
$ws.Range("D7").Value = "'0.3812"
$ws.Range("E7").Value = "'  +1.36%  "
$ws.Range("D8").Value = "'0.3625"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("D9").Value = "'51.52"
$ws.Range("E9").Value = "'  +0.17%  "
$ws.Range("D10").Value = "'1.261"
$ws.Range("E10").Value = "'  +3.35%  "
$ws.Range("D11").Value = "'0.08229"
$ws.Range("E11").Value = "'  +0.92%  "
$ws.Range("E12").Value = "'  +0.11%  "
$ws.Range("D13").Value = "'22.73"
$ws.Range("E13").Value = "'  +2.22%  "
$ws.Range("E14").Value = "'  +1.26%  "
$ws.Range("D15").Value = "'7.470"
$ws.Range("E15").Value = "'  +2.43%  "
$ws.Range("E16").Value = "'  +0.42%  "
$ws.Range("D17").Value = "'1.658.69"
$ws.Range("E17").Value = "'  +3.36%  "
$ws.Range("D18").Value = "'97.97"
$ws.Range("E18").Value = "'  +3.87%  "
$ws.Range("E19").Value = "'  +0.78%  "
$ws.Range("D20").Value = "'6.801"
$ws.Range("E20").Value = "'  +3.99%  "
$ws.Range("D21").Value = "'17.79"
$ws.Range("E21").Value = "'  +1.37%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  +0.07%  "
$ws.Range("D23").Value = "'12.80"
$ws.Range("E23").Value = "'  +2.41%  "
$ws.Range("D24").Value = "'23.785.67"
$ws.Range("E24").Value = "'  +1.90%  "
$ws.Range("D25").Value = "'2.560"
$ws.Range("E25").Value = "'  +3.41%  "
$ws.Range("E26").Value = "'  +0.46%  "
$ws.Range("D27").Value = "'21.32"
$ws.Range("E27").Value = "'  +0.97%  "
$ws.Range("D28").Value = "'151.26"
$ws.Range("E28").Value = "'  +0.73%  "
$ws.Range("D29").Value = "'5.238"
$ws.Range("E29").Value = "'  -0.71%  "
$ws.Range("D30").Value = "'134.60"
$ws.Range("E30").Value = "'  +1.35%  "
$ws.Range("D31").Value = "'1.839.08"
$ws.Range("E31").Value = "'  +2.51%  "
$ws.Range("D32").Value = "'6.953"
$ws.Range("E32").Value = "'  +3.73%  "
$ws.Range("E33").Value = "'  +2.10%  "
$ws.Range("D34").Value = "'1.085"
$ws.Range("E34").Value = "'  +3.02%  "
$ws.Range("D35").Value = "'11.90"
$ws.Range("E35").Value = "'  +6.35%  "
$ws.Range("E36").Value = "'  +2.60%  "
$ws.Range("D37").Value = "'6.167"
$ws.Range("E37").Value = "'  +3.07%  "
$ws.Range("D38").Value = "'0.2523"
$ws.Range("E38").Value = "'  +1.38%  "
$ws.Range("D39").Value = "'0.08824"
$ws.Range("E39").Value = "'  +0.67%  "
$ws.Range("D40").Value = "'0.07136"
$ws.Range("E40").Value = "'  +0.22%  "
$ws.Range("D41").Value = "'13.38"
$ws.Range("E41").Value = "'  +11.24%  "
$ws.Range("D42").Value = "'0.7080"
$ws.Range("E42").Value = "'  +1.41%  "
$ws.Range("D43").Value = "'1.346"
$ws.Range("E43").Value = "'  +1.09%  "
$ws.Range("D44").Value = "'15.97"
$ws.Range("E44").Value = "'  +0.81%  "
$ws.Range("D45").Value = "'0.6571"
$ws.Range("E45").Value = "'  +1.89%  "
$ws.Range("D46").Value = "'2.333"
$ws.Range("E46").Value = "'  +2.60%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "'  +0.14%  "
$ws.Range("D48").Value = "'3.959"
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E49").Value = "'  +0.10%  "
$ws.Range("D50").Value = "'128.48"
$ws.Range("E50").Value = "'  +0.99%  "
$ws.Range("D51").Value = "'1.196"
$ws.Range("E51").Value = "'  +0.47%  "
